$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Update existing metric values in row 2
$ws.Range("B2").Value = 0.08969192748125149
$ws.Range("C2").Value = 0.9987856266669743
$ws.Range("D2").Value = 0.2402316749144211

# Update the model description text (includes a line break + indentation)
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 AdaBoostRegressor(learning_rate=0.5, n_estimators=100))])"

# Add the new Elapsed Time / CPU values
$ws.Range("G2").Value = 0.1289622459000384
$ws.Range("H2").Value = 0.991
